$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A27").Value = 52
$ws.Range("B27").Value = "okay added"
$ws.Range("C27").Value = "riya-morankar"
$ws.Range("D27").Value = "N/A"
$ws.Range("E27").Value = "edit1 to main"
$ws.Range("F27").NumberFormat = "@"
$ws.Range("F27").Value = "2025-06-20"
